$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 42612.884166666663
$ws.Range("B8").Value = -12
$ws.Range("C8").Value = 49
$ws.Range("D8").Value = 50
$ws.Range("E8").Value = 22
$ws.Range("F8").Value = 77
$ws.Range("G8").Value = 11832
$ws.Range("H8").Value = 14416
$ws.Range("I8").Value = 1563
$ws.Range("J8").Value = 142
$ws.Range("K8").Value = 143
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = 7
$ws.Range("N8").Value = "Bag"

$ws.Range("A9").Value = 42612.890694444446
$ws.Range("B9").Value = -10
$ws.Range("C9").Value = 50
$ws.Range("D9").Value = 48
$ws.Range("E9").Value = 22
$ws.Range("F9").Value = 77
$ws.Range("G9").Value = 9089
$ws.Range("H9").Value = 14398
$ws.Range("I9").Value = 1565
$ws.Range("J9").Value = 152
$ws.Range("K9").Value = 146
$ws.Range("L9").Value = 2
$ws.Range("M9").Value = 7
$ws.Range("N9").Value = "Bag"
